# 9.3.1.xlsx update: add a new "2021" year column (O) to the small-scale
# industries table, update a couple of existing year values, and move the
# active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting for the new column O by copying it from column N,
#     which holds the equivalent "last year" formatting for each row. This
#     keeps border/font/number-format/alignment identical without
#     introducing brand-new style definitions.

# Row 3 (thin separator row above the header) - empty cell, border only.
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)

# Row 4 (year header row) - new header value 2021.
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

# Row 5 (data row) - new data value for 2021.
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 4.1

# --- Update existing data values in row 5.
$ws.Range("L5").Value = 1.6
$ws.Range("N5").Value = 3.1

# --- Move the active selection (matches the saved view state).
[void]$ws.Range("P4").Select()

Write-Host "9.3.1.xlsx updated: column O (2021) added, L5/N5 refreshed, selection moved to P4"
